## Loan RBI, Variable Instalments
## On the "Repayment Schedule" sheet a new (blank) column is inserted
## before the existing "Late" column (column N), pushing "Late" into O
## and "Outstanding" into Q. The "Repayment Schedule" sheet also becomes
## the active/selected sheet (it was "Input" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Capture the width of the column immediately to the left (M) so the
# freshly inserted column inherits the same look & feel.
$existingWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column at N - this shifts the old N ("Late") to O
# and the old O/P ("Outstanding" / values) to P/Q respectively.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $existingWidth

# Make "Repayment Schedule" the active sheet/tab (it moves the
# tabSelected flag off of "Input" and onto this sheet automatically).
$ws.Activate()
$ws.Range("R8").Select()
